# Pavel - new user for linking test
# Add a new row to the "Users" sheet for the Linking_AutoUser test account.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Row 44 is an existing data row that carries the plain "Normal" style
# (thin border, no special fill/font) uniformly across columns A:G - copy
# its formatting down onto the new row 52 so the new row matches its
# neighbours visually.
$ws.Range("A44:G44").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new user row. C52/D52 are intentionally left blank, matching
# the pattern used by the other "test-only" user rows (e.g. row 50, 51).
$ws.Range("A52").Value = "Linking_AutoUser"
$ws.Range("B52").Value = "Password1"
$ws.Range("E52").Value = "Default user for Linking tests"
$ws.Range("F52").Value = "N"
$ws.Range("G52").Value = "linking.autouser@mailinator.com"

# Restore the view/selection state that was active when the row was added.
$ws.Activate()
$ws.Range("L33").Select()
